$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C21").Value = 52
$ws.Range("D21").Value = 46
$ws.Range("E21").Value = 6
$ws.Range("F21").Value = 13.18051575931232
